$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision change on the existing A25 timestamp
$ws.Cells.Item(25, 1).Value = 44338.77958125116

# Append the new data row (row 26) retrieved in this run
$ws.Cells.Item(26, 1).Value = 44339.77832643608
$ws.Cells.Item(26, 2).Value = 74327
$ws.Cells.Item(26, 3).Value = 62568
$ws.Cells.Item(26, 4).Value = 3340
$ws.Cells.Item(26, 5).Value = 2091
$ws.Cells.Item(26, 6).Value = 1475
$ws.Cells.Item(26, 7).Value = 19302
$ws.Cells.Item(26, 8).Value = 1416
$ws.Cells.Item(26, 9).Value = 834
$ws.Cells.Item(26, 10).Value = 213
